$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = 41599
$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("C3").Select()
